$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'264.85"
$c.Style = 'Normal'
$c = $ws.Range('G2')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G3')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D4')
$c.Value = "'6.225"
$c.Style = 'Normal'
$c = $ws.Range('G4')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G5')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.Value = "'3.552"
$c.Style = 'Normal'
$c = $ws.Range('G6')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D7')
$c.Value = "'6.733"
$c.Style = 'Normal'
$c = $ws.Range('G7')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G8')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D9')
$c.Value = "'0.8145"
$c.Style = 'Normal'
$c = $ws.Range('G9')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.Value = "'0.1594"
$c.Style = 'Normal'
$c = $ws.Range('G10')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.Value = "'0.08233"
$c.Style = 'Normal'
$c = $ws.Range('G11')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.Value = "'0.03385"
$c.Style = 'Normal'
$c = $ws.Range('G12')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.Value = "'0.03151"
$c.Style = 'Normal'
$c = $ws.Range('G13')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D14')
$c.Value = "'0.09269"
$c.Style = 'Normal'
$c = $ws.Range('G14')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.Value = "'3.897"
$c.Style = 'Normal'
$c = $ws.Range('G15')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.Value = "'0.001703"
$c.Style = 'Normal'
$c = $ws.Range('G16')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.Value = "'0.04843"
$c.Style = 'Normal'
$c = $ws.Range('G17')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.Value = "'0.0006270"
$c.Style = 'Normal'
$c = $ws.Range('G18')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.Value = "'0.006220"
$c.Style = 'Normal'
$c = $ws.Range('G19')
$c.Value = "'14"
$c.Style = 'Normal'
$ws.Range('B20').Value = 'HotbitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$c = $ws.Range('D20')
$c.Value = "'0.006215"
$c.Style = 'Normal'
$ws.Range('E20').Value = '19HotbitTokenHTB'
$c = $ws.Range('G20')
$c.Value = "'14"
$c.Style = 'Normal'
$ws.Range('B21').Value = 'BitKan'
$ws.Range('C21').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$c = $ws.Range('D21')
$c.Value = "'0.001101"
$c.Style = 'Normal'
$ws.Range('E21').Value = '20BitKanKAN'
$c = $ws.Range('G21')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D22')
$c.Value = "'0.0001503"
$c.Style = 'Normal'
$c = $ws.Range('G22')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D23')
$c.Value = "'3.697"
$c.Style = 'Normal'
$c = $ws.Range('G23')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.Value = "'2.261"
$c.Style = 'Normal'
$c = $ws.Range('G24')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G25')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G26')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.Value = "'0.0002691"
$c.Style = 'Normal'
$c = $ws.Range('G27')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G28')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G29')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G30')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G31')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G32')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G33')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G34')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G35')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G36')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G37')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G38')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G39')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D40')
$c.Value = "'0.04608"
$c.Style = 'Normal'
$c = $ws.Range('G40')
$c.Value = "'14"
$c.Style = 'Normal'
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$c = $ws.Range('D41')
$c.Value = "'0.1128"
$c.Style = 'Normal'
$ws.Range('E41').Value = '40BKEXTokenBKK'
$c = $ws.Range('G41')
$c.Value = "'14"
$c.Style = 'Normal'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$c = $ws.Range('D42')
$c.Value = "'0.003137"
$c.Style = 'Normal'
$ws.Range('E42').Value = '41CEJICEJI'
$c = $ws.Range('G42')
$c.Value = "'14"
$c.Style = 'Normal'
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$c = $ws.Range('D43')
$c.Value = "'0.007266"
$c.Style = 'Normal'
$ws.Range('E43').Value = '42KickTokenKICK'
$c = $ws.Range('G43')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G44')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D45')
$c.Value = "'0.00006178"
$c.Style = 'Normal'
$c = $ws.Range('G45')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G46')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.Value = "'0.7525"
$c.Style = 'Normal'
$c = $ws.Range('G47')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D48')
$c.Value = "'0.1725"
$c.Style = 'Normal'
$c = $ws.Range('G48')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.Value = "'0.00002108"
$c.Style = 'Normal'
$c = $ws.Range('G49')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.Value = "'0.01244"
$c.Style = 'Normal'
$c = $ws.Range('G50')
$c.Value = "'14"
$c.Style = 'Normal'
$c = $ws.Range('G51')
$c.Value = "'14"
$c.Style = 'Normal'
